$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "'257.45"
$ws.Range("E2").Value2 = "'-0.26%"
$ws.Range("D3").Value2 = "'27.45"
$ws.Range("E3").Value2 = "'-1.26%"
$ws.Range("D4").Value2 = "'4.581"
$ws.Range("E4").Value2 = "'-12.28%"
$ws.Range("D5").Value2 = "'0.05885"
$ws.Range("E5").Value2 = "'-1.07%"
$ws.Range("D6").Value2 = "'6.630"
$ws.Range("E6").Value2 = "'-1.07%"
$ws.Range("D7").Value2 = "'0.8578"
$ws.Range("E7").Value2 = "'-1.61%"
$ws.Range("E8").Value2 = "'-10.94%"
$ws.Range("E9").Value2 = "'-1.12%"
$ws.Range("D10").Value2 = "'0.03583"
$ws.Range("E10").Value2 = "'-1.21%"
$ws.Range("D11").Value2 = "'0.07081"
$ws.Range("E11").Value2 = "'-2.24%"
$ws.Range("E12").Value2 = "'0.06%"
$ws.Range("D13").Value2 = "'0.09177"
$ws.Range("E13").Value2 = "'-0.62%"
$ws.Range("D14").Value2 = "'0.001543"
$ws.Range("E14").Value2 = "'-1.24%"
$ws.Range("D15").Value2 = "'0.0006067"
$ws.Range("E15").Value2 = "'-0.05%"
$ws.Range("D16").Value2 = "'0.006069"
$ws.Range("E16").Value2 = "'1.65%"
$ws.Range("D17").Value2 = "'3.517"
$ws.Range("E17").Value2 = "'0.66%"
$ws.Range("D18").Value2 = "'3.197"
$ws.Range("E18").Value2 = "'-2.16%"
$ws.Range("E19").Value2 = "'-0.27%"
$ws.Range("D20").Value2 = "'0.3108"
$ws.Range("E20").Value2 = "'-1.26%"
$ws.Range("E21").Value2 = "'-1.06%"
$ws.Range("D22").Value2 = "'3.844"
$ws.Range("E22").Value2 = "'8.85%"
$ws.Range("D23").Value2 = "'0.04214"
$ws.Range("E23").Value2 = "'0.68%"
$ws.Range("D24").Value2 = "'0.001218"
$ws.Range("E24").Value2 = "'-0.19%"
$ws.Range("D25").Value2 = "'0.004302"
$ws.Range("E25").Value2 = "'-5.51%"
$ws.Range("D26").Value2 = "'0.0001199"
$ws.Range("E26").Value2 = "'-0.22%"
$ws.Range("E27").Value2 = "'-22.09%"
$ws.Range("D40").Value2 = "'0.03834"
$ws.Range("E40").Value2 = "'0.29%"
$ws.Range("B41").Value2 = 'BKEXToken'
$ws.Range("C41").Value2 = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("D41").Value2 = "'0.1103"
$ws.Range("E41").Value2 = "'-0.66%"
$ws.Range("B42").Value2 = 'KickToken'
$ws.Range("C42").Value2 = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range("D42").Value2 = "'0.003963"
$ws.Range("E42").Value2 = "'-28.99%"
$ws.Range("D43").Value2 = "'0.002399"
$ws.Range("E43").Value2 = "'0.66%"
$ws.Range("D44").Value2 = "'0.01179"
$ws.Range("E44").Value2 = "'18.71%"
$ws.Range("D45").Value2 = "'0.00005468"
$ws.Range("E45").Value2 = "'0.70%"
$ws.Range("E46").Value2 = "'-0.03%"
$ws.Range("D47").Value2 = "'0.05997"
$ws.Range("E47").Value2 = "'-45.03%"
$ws.Range("D48").Value2 = "'0.1346"
$ws.Range("E48").Value2 = "'6,183.07%"
$ws.Range("D49").Value2 = "'0.00002099"
$ws.Range("E49").Value2 = "'-0.03%"
$ws.Range("D50").Value2 = "'0.0001999"
$ws.Range("E50").Value2 = "'-0.03%"
